$p = $ppt.ActivePresentation

# --- Create the new slide 6 as a duplicate of slide 5 (the most recent "logo" slide). ---
$src = $p.Slides.Item(5)
$dup = $src.Duplicate()
$s6 = $p.Slides.Item($p.Slides.Count)

# --- Ungroup the big "Group 2" container that held all seven logo shapes so the ---
# --- two plain ovals + five sub-groups become direct children of the slide, then ---
# --- drop the Heart + black-dot shapes that aren't part of this logo variant.   ---
$container = $s6.Shapes.Item(1)
$container.Ungroup() | Out-Null

for ($i = $s6.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s6.Shapes.Item($i)
    if ($sh.Name -eq "Heart 11" -or $sh.Name -eq "Oval 12") {
        $sh.Delete()
    }
}

# --- Spread the seven logo pieces out along a single horizontal row. ---
function Set-ShapePos($shape, [double]$left, [double]$top) {
    $shape.Left = $left
    $shape.Top = $top
}

for ($i = 1; $i -le $s6.Shapes.Count; $i++) {
    $sh = $s6.Shapes.Item($i)
    switch ($sh.Id) {
        5  { Set-ShapePos $sh 265.93544007086615 92.90322494645669 }
        25 { Set-ShapePos $sh 332.516067592126   92.90322494645669 }
        26 { Set-ShapePos $sh 399.09669501338584 92.90322494645669 }
        27 { Set-ShapePos $sh 465.6773224346457  92.90322494645669 }
        2  { Set-ShapePos $sh 532.6452941905511  92.90322494645669 }
        29 { Set-ShapePos $sh 599.6132507464567  92.90322494645669 }
    }
}

# --- Re-flow the caption text box above the new single-row layout. ---
for ($i = 1; $i -le $s6.Shapes.Count; $i++) {
    $sh = $s6.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 30") {
        $sh.Left = 199.3548050496063
        $sh.Top = 173.8063736527559
        $sh.Width = 464.12944031889765
        $sh.Height = 60.58590511181102
    }
}
